$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting existing D:K to E:L
$ws.Columns("D").Insert()

# Copy number formats/styles from the (now-shifted) old D column (now in E) into new D
$ws.Range("E5:E102").Copy() | Out-Null
$ws.Range("D5:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate new column D with the new quarter (period ending 2018-11-03) data
$ws.Range("D7").Value = 43407
$ws.Range("D8").Value = 9825800
$ws.Range("D9").Value = 6995200
$ws.Range("D10").Value = 2830600
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 36100
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 8788500
$ws.Range("D18").Value = 1037300
$ws.Range("D20").Value = 25700
$ws.Range("D21").Value = 1267900
$ws.Range("D22").Value = 16500
$ws.Range("D23").Value = 1046500
$ws.Range("D24").Value = 284300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 762300
$ws.Range("D27").Value = 762300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -25700
$ws.Range("D33").Value = 762300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 762300
$ws.Range("D38").Value = 43407
$ws.Range("D41").Value = 2711800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 517400
$ws.Range("D44").Value = 5543400
$ws.Range("D45").Value = 544400
$ws.Range("D46").Value = 9317000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 5165900
$ws.Range("D49").Value = 97300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 445000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 15025200
$ws.Range("D57").Value = 3340600
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 2673200
$ws.Range("D60").Value = 6013800
$ws.Range("D61").Value = 2232900
$ws.Range("D62").Value = 1521700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 9768400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 4615600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 5256900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43407
$ws.Range("D81").Value = 762300
$ws.Range("D83").Value = 204900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 917100
$ws.Range("D91").Value = -299100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -298200
$ws.Range("D96").Value = -241400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -769900
$ws.Range("D101").Value = -10000
$ws.Range("D102").Value = -161000
